# Add a new worksheet "ODI Batting Extra" as the last sheet and populate it
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $lastSheet)
$new.Name = "ODI Batting Extra"

# Header row - reuse the bold/bordered header style already used on the other
# sheets by copying it, then overwrite the text.
$ws1.Range("A1:D1").Copy($new.Range("A1:F1"))

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $headers.Length; $c++) {
    $new.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Data rows - MATCH_CODE (A), NUM_4 (C), NUM_6 (D) and PERCENT_RUNS_OF_TOTAL (E)
# are text values even though numeric looking, so force text format before
# writing them. BATTING_POSITION (B) stays numeric. MAN_OF_MATCH (F) is text.
$new.Range("A2:A7").NumberFormat = "@"
$new.Range("C2:E7").NumberFormat = "@"
$new.Range("F2:F7").NumberFormat = "@"

$data = @(
    @("4273", 2, "5", "0", "14.29%", "NO"),
    @("4274", 2, "4", "0", "6.69%", "NO"),
    @("4275", 2, "0", "0", "1.08%", "NO"),
    @("4276", $null, $null, $null, $null, "NO"),
    @("4277", 1, "3", "0", "16.29%", "NO"),
    @("4690", $null, $null, $null, $null, "NO")
)

$r = 2
foreach ($row in $data) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $val = $row[$c - 1]
        if ($null -ne $val) {
            $new.Cells.Item($r, $c).Value = $val
        }
    }
    $r++
}

# Restore the originally active sheet/selection so the workbook view state
# is left as it was before the edit.
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
